$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("NV-6 Lâm Thị Mỹ Hằng")
$ws.Range("B11").Value = 2515000
$ws.Range("B16").Value = 5750000
$ws.Range("B19").Value = 2100000
$ws.Range("B27").Value = 3950000
$ws.Range("B29").Value = 4600000
$ws.Range("B35").Value = 6794047.619047619
$ws.Range("B36").Value = 14269047.61904762
$ws.Range("B37").Value = 11219047.61904762
$ws.Range("B38").Value = 32282142.85714286

$ws = $wb.Worksheets.Item("NV-7 Phạm Thanh Hoàng")
$ws.Range("B11").Value = 1980000
$ws.Range("B20").Value = 2100000
$ws.Range("B30").Value = 1800000
$ws.Range("B36").Value = 3166785.714285716
$ws.Range("B37").Value = 3760714.285714285
$ws.Range("B38").Value = 4389285.714285715
$ws.Range("B39").Value = 11316785.71428572

$ws = $wb.Worksheets.Item("NV-23 Lê Hoàng Thanh")
$ws.Range("B20").Value = -1000000
$ws.Range("B32").Value = 2928571.428571429
$ws.Range("B34").Value = 3628571.428571429

$ws = $wb.Worksheets.Item("NV-40 Sang sang")
$ws.Range("B21").Value = -2500000
$ws.Range("B33").Value = 1228571.428571429
$ws.Range("B35").Value = 1228571.428571429

$ws = $wb.Worksheets.Item("NV-10 Lê Đình Hậu")
$ws.Range("B27").Value = 6660000
$ws.Range("B30").Value = 0
$ws.Range("B34").Value = -11100000
$ws.Range("B37").Value = -2857.142857141793
$ws.Range("B38").Value = 7807142.857142858

$ws = $wb.Worksheets.Item("NV-16 Kha Như Huỳnh")
$ws.Range("B27").Value = 400000
$ws.Range("B29").Value = 500000
$ws.Range("B34").Value = 348571.4285714282
$ws.Range("B35").Value = 348571.4285714282

$ws = $wb.Worksheets.Item("NV-11 Đỗ Thị Huyền Trân")
$ws.Range("B5").Value = 2992000
$ws.Range("B11").Value = 0
$ws.Range("B35").Value = 7880214.285714285
$ws.Range("B38").Value = 7880214.285714285

$ws = $wb.Worksheets.Item("NV-26 Trần Khánh Hiệp")
$ws.Range("B29").Value = 100000
$ws.Range("B34").Value = -3671428.571428571
$ws.Range("B35").Value = -3671428.571428571

$ws = $wb.Worksheets.Item("NV-29 Lâm Hoàng Phú")
$ws.Range("B10").Value = 400000
$ws.Range("B34").Value = 489285.7142857141
$ws.Range("B37").Value = 589285.7142857141

$ws = $wb.Worksheets.Item("NV-5 Nguyễn Hoàng Yến Quyên")
$ws.Range("B7").Value = 5122000
$ws.Range("B17").Value = 2800000
$ws.Range("B27").Value = 3820000
$ws.Range("B32").Value = 12257000
$ws.Range("B33").Value = 3657142.857142857
$ws.Range("B34").Value = 4677142.857142857
$ws.Range("B35").Value = 20591285.71428572
